$wb = $excel.ActiveWorkbook

# Replace the status text "Ready for handoff" with "In Translation"
# everywhere it appears (Overview sheet columns E/F, zh-cn/de-de sheet column C).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# Narrow the "Status"-related columns (target OOXML width ~13.41 chars).
# NOTE: the host quantizes ColumnWidth to 1/6-character increments, so
# 12.5 is the input that lands closest to the recorded 13.4101845877511.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
